$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.617.28"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "1.962.89"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.618"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.73"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.379"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0805"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.93%  "
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D14").Value = "2.250.04"
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").Value = "1.976.80"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "36.536.07"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").Value = "0.0₃0857"
$ws.Range("E20").Value = "  -2.64%  "
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.139"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("E30").Value = "  +1.28%  "
$ws.Range("E31").Value = "  -3.29%  "
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("E33").Value = "  -3.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.32"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("E36").Value = "  +2.04%  "
$ws.Range("E37").Value = "  +10.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -10.44%  "
$ws.Range("E40").Value = "  -2.13%  "
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.28%  "
$ws.Range("D45").Value = "1.365.45"
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").Value = "2.141.02"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.34%  "
